$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3999.5
$ws.Range("J100").Value = 3999.5
$ws.Range("L100").Value = 3999.5
$ws.Range("N100").Value = -5081.5
$ws.Range("H118").Value = 1003
$ws.Range("I118").Value = 803.6
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 2410.8
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = -753.8000000000002
$ws.Range("N118").Value = -9314
$ws.Range("H132").Value = 1299
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H137").Value = 1481.9048
$ws.Range("I137").Value = 759.7143
$ws.Range("J137").Value = 2926.2856
$ws.Range("K137").Value = 2279.1429
$ws.Range("L137").Value = 8778.856800000001
$ws.Range("M137").Value = 270.8571000000002
$ws.Range("N137").Value = -13878.8568
$ws.Range("H138").Value = 4362.647
$ws.Range("I138").Value = 2476.9
$ws.Range("J138").Value = 5148.375
$ws.Range("K138").Value = 7430.700000000001
$ws.Range("L138").Value = 15445.125
$ws.Range("M138").Value = -2290.700000000001
$ws.Range("N138").Value = -25725.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 29450
$ws.Range("J24").Value = 29450
$ws.Range("L24").Value = 29450
$ws.Range("N24").Value = -30198
$ws.Range("H61").Value = 812
$ws.Range("I61").Value = 812
$ws.Range("K61").Value = 812
$ws.Range("M61").Value = -600
$ws.Range("H97").Value = 1762.5
$ws.Range("I97").Value = 1720.8
$ws.Range("K97").Value = 1720.8
$ws.Range("M97").Value = -1224.8
$ws.Range("H100").Value = 29450
$ws.Range("J100").Value = 29450
$ws.Range("L100").Value = 29450
$ws.Range("N100").Value = -31614
$ws.Range("H112").Value = 28283.334
$ws.Range("J112").Value = 28283.334
$ws.Range("L112").Value = 28283.334
$ws.Range("N112").Value = -31237.334
$ws.Range("H114").Value = 38000
$ws.Range("J114").Value = 38000
$ws.Range("L114").Value = 38000
$ws.Range("N114").Value = -46678
$ws.Range("H119").Value = 89999
$ws.Range("J119").Value = 89999
$ws.Range("L119").Value = 89999
$ws.Range("N119").Value = -99675
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820
$ws.Range("H132").Value = 3235.889
$ws.Range("I132").Value = 3235.889
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9707.667000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7177.667000000001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 812
$ws.Range("I136").Value = 812
$ws.Range("K136").Value = 2436
$ws.Range("M136").Value = 114

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H19").Value = 100005
$ws.Range("J19").Value = 100005
$ws.Range("L19").Value = 100005
$ws.Range("N19").Value = -100351
$ws.Range("H111").Value = 50696
$ws.Range("J111").Value = 50696
$ws.Range("L111").Value = 50696
$ws.Range("N111").Value = -58876
$ws.Range("H134").Value = 2953.7144
$ws.Range("I134").Value = 3112.5
$ws.Range("J134").Value = 2001
$ws.Range("K134").Value = 9337.5
$ws.Range("L134").Value = 6003
$ws.Range("M134").Value = -6802.5
$ws.Range("N134").Value = -11073

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 516.4
$ws.Range("I22").Value = 486
$ws.Range("K22").Value = 486
$ws.Range("M22").Value = -136
$ws.Range("H31").Value = 1487.1428
$ws.Range("I31").Value = 1282
$ws.Range("K31").Value = 1282
$ws.Range("M31").Value = -987
$ws.Range("H34").Value = 1487.1428
$ws.Range("I34").Value = 1282
$ws.Range("K34").Value = 1282
$ws.Range("M34").Value = -1080
$ws.Range("H58").Value = 2024.75
$ws.Range("I58").Value = 1899.6666
$ws.Range("K58").Value = 1899.6666
$ws.Range("M58").Value = -1696.6666
$ws.Range("H92").Value = 28985.666
$ws.Range("J92").Value = 28985.666
$ws.Range("L92").Value = 28985.666
$ws.Range("N92").Value = -33977.666
$ws.Range("H105").Value = 3041.0625
$ws.Range("I105").Value = 2377.6
$ws.Range("J105").Value = 3342.6365
$ws.Range("K105").Value = 2377.6
$ws.Range("L105").Value = 3342.6365
$ws.Range("M105").Value = -630.5999999999999
$ws.Range("N105").Value = -6836.636500000001
$ws.Range("H132").Value = 4254.7144
$ws.Range("I132").Value = 5848
$ws.Range("J132").Value = 2130.3333
$ws.Range("K132").Value = 17544
$ws.Range("L132").Value = 6390.999899999999
$ws.Range("M132").Value = -15014
$ws.Range("N132").Value = -11450.9999
$ws.Range("H136").Value = 2024.75
$ws.Range("I136").Value = 1899.6666
$ws.Range("K136").Value = 5698.9998
$ws.Range("M136").Value = -3148.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1792.8
$ws.Range("I51").Value = 449.5
$ws.Range("J51").Value = 2128.625
$ws.Range("K51").Value = 1348.5
$ws.Range("L51").Value = 6385.875
$ws.Range("M51").Value = -888.5
$ws.Range("N51").Value = -7305.875
$ws.Range("H75").Value = 6619.5
$ws.Range("J75").Value = 9303.5
$ws.Range("L75").Value = 27910.5
$ws.Range("N75").Value = -29906.5
$ws.Range("H78").Value = 6619.5
$ws.Range("J78").Value = 9303.5
$ws.Range("L78").Value = 83731.5
$ws.Range("N78").Value = -93715.5
$ws.Range("H97").Value = 435.2
$ws.Range("J97").Value = 425.8
$ws.Range("L97").Value = 1277.4
$ws.Range("N97").Value = -2269.4
$ws.Range("H128").Value = 617416.3
$ws.Range("I128").Value = 617416.3
$ws.Range("K128").Value = 1852248.9
$ws.Range("M128").Value = -1847268.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 12603.143
$ws.Range("J98").Value = 12603.143
$ws.Range("L98").Value = 12603.143
$ws.Range("N98").Value = -18593.143
$ws.Range("H102").Value = 3368.6667
$ws.Range("I102").Value = 3368.6667
$ws.Range("K102").Value = 3368.6667
$ws.Range("M102").Value = -1746.6667
$ws.Range("H122").Value = 2167.8572
$ws.Range("I122").Value = 1862.5
$ws.Range("K122").Value = 5587.5
$ws.Range("M122").Value = -3137.5
$ws.Range("H126").Value = 3447
$ws.Range("I126").Value = 4444
$ws.Range("J126").Value = 2450
$ws.Range("K126").Value = 13332
$ws.Range("L126").Value = 7350
$ws.Range("M126").Value = -10862
$ws.Range("N126").Value = -12290
$ws.Range("H132").Value = 2999.5
$ws.Range("I132").Value = 2999.5
$ws.Range("K132").Value = 8998.5
$ws.Range("M132").Value = -6468.5
$ws.Range("H140").Value = 143958.33
$ws.Range("J140").Value = 143958.33
$ws.Range("L140").Value = 143958.33
$ws.Range("N140").Value = -154318.33

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 19000
$ws.Range("J69").Value = 19000
$ws.Range("L69").Value = 19000
$ws.Range("N69").Value = -20622
$ws.Range("H72").Value = 19000
$ws.Range("J72").Value = 19000
$ws.Range("L72").Value = 57000
$ws.Range("N72").Value = -65112
$ws.Range("H93").Value = 3370
$ws.Range("I93").Value = 3740
$ws.Range("K93").Value = 3740
$ws.Range("M93").Value = -2492
$ws.Range("H119").Value = 150000
$ws.Range("J119").Value = 150000
$ws.Range("L119").Value = 150000
$ws.Range("N119").Value = -159676
$ws.Range("H124").Value = 12000
$ws.Range("J124").Value = 12000
$ws.Range("L124").Value = 12000
$ws.Range("N124").Value = -21820
$ws.Range("H127").Value = 49500
$ws.Range("J127").Value = 49500
$ws.Range("L127").Value = 49500
$ws.Range("N127").Value = -59420
$ws.Range("H132").Value = 5703.6665
$ws.Range("I132").Value = 6046.5
$ws.Range("K132").Value = 18139.5
$ws.Range("M132").Value = -15609.5
$ws.Range("H136").Value = 1954.8695
$ws.Range("I136").Value = 886.55554
$ws.Range("J136").Value = 5800.8
$ws.Range("K136").Value = 2659.66662
$ws.Range("L136").Value = 17402.4
$ws.Range("M136").Value = -109.66662
$ws.Range("N136").Value = -22502.4
$ws.Range("H138").Value = 130000
$ws.Range("J138").Value = 130000
$ws.Range("L138").Value = 130000
$ws.Range("N138").Value = -140280

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19983.2
$ws.Range("J41").Value = 19987
$ws.Range("L41").Value = 19987
$ws.Range("N41").Value = -20767
$ws.Range("H55").Value = 9400.333000000001
$ws.Range("I55").Value = 5011.5
$ws.Range("K55").Value = 5011.5
$ws.Range("M55").Value = -4734.5
$ws.Range("H61").Value = 9999
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 9999
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 9999
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -10583
$ws.Range("H141").Value = 219331.83
$ws.Range("J141").Value = 220398.4
$ws.Range("L141").Value = 220398.4
$ws.Range("N141").Value = -230758.4
